# Updates the cryptos price/volume snapshot to the values scraped on the
# latest GitHub Actions run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 2-43 (excluding the mCoin/RenderToken/MXToken block) -------------
$ws.Range("D2").Value = "28.131.81"
$ws.Range("E2").Value = "  +2.33%  "
$ws.Range("D3").Value = "1.654.48"
$ws.Range("E3").Value = "  +2.35%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'214.08"
$ws.Range("E5").Value = "  +1.56%  "
$ws.Range("D6").Value = "'0.529"
$ws.Range("E6").Value = "  +0.61%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("D8").Value = "'23.62"
$ws.Range("E8").Value = "  +3.77%  "
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("D10").Value = "'0.0615"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").Value = "'0.0874"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").Value = "1.886.25"
$ws.Range("E12").Value = "  +2.18%  "
$ws.Range("D13").Value = "1.655.26"
$ws.Range("E13").Value = "  +2.45%  "
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").Value = "'0.568"
$ws.Range("E15").Value = "  +3.50%  "
$ws.Range("D16").Value = "'65.81"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").Value = "28.113.60"
$ws.Range("E17").Value = "  +2.34%  "
$ws.Range("D18").Value = "'233.84"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").Value = "'7.72"
$ws.Range("E19").Value = "  +2.66%  "
$ws.Range("D20").Value = "0.0₃0725"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'10.70"
$ws.Range("E22").Value = "  +5.60%  "
$ws.Range("D23").Value = "'4.42"
$ws.Range("E23").Value = "  +3.10%  "
$ws.Range("D24").Value = "'2.16"
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("D25").Value = "'152.35"
$ws.Range("E25").Value = "  +1.20%  "
$ws.Range("D26").Value = "'6.94"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("D27").Value = "'15.82"
$ws.Range("E27").Value = "  +1.94%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("D31").Value = "'0.0485"
$ws.Range("E31").Value = "  +0.51%  "
$ws.Range("E32").Value = "  +2.64%  "
$ws.Range("D33").Value = "1.453.71"
$ws.Range("E33").Value = "  -0.70%  "
$ws.Range("E34").Value = "  +0.45%  "
$ws.Range("E35").Value = "  +2.68%  "
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("D37").Value = "'0.893"
$ws.Range("E37").Value = "  +4.05%  "
$ws.Range("D38").Value = "'0.0170"
$ws.Range("E38").Value = "  +1.53%  "
$ws.Range("D39").Value = "'0.561"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D40").Value = "'0.926"
$ws.Range("E40").Value = "  -1.72%  "
$ws.Range("D41").Value = "'69.57"
$ws.Range("E41").Value = "  +2.33%  "
$ws.Range("E42").Value = "  +3.44%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.19%  "

# --- rows 44-46: coins re-ranked (MXToken, mCoin, RenderToken) -------------
$ws.Range("B44").Value = "MXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D44").Value = "'2.30"
$ws.Range("E44").Value = "  +3.99%  "

$ws.Range("B45").Value = "mCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/fzVgyjBcRc9+mcoin-mcoin"
$ws.Range("D45").Value = "'2.46"
$ws.Range("E45").Value = "  -0.54%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'1.82"
$ws.Range("E46").Value = "  +6.21%  "

# --- rows 47-51 --------------------------------------------------------------
$ws.Range("D47").Value = "'5.42"
$ws.Range("E47").Value = "  +3.06%  "
$ws.Range("D48").Value = "1.797.00"
$ws.Range("E48").Value = "  +2.11%  "
$ws.Range("D49").Value = "'89.07"
$ws.Range("E49").Value = "  +2.88%  "
$ws.Range("D50").Value = "'0.102"
$ws.Range("E50").Value = "  +0.86%  "
$ws.Range("D51").Value = "'7.74"
$ws.Range("E51").Value = "  +1.48%  "
